$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loss-of-sale rows to append after the existing data (row 43 was the last
# existing row). Columns: #, Date, Customer Name, Contact, Function Date,
# Staff, Status, Category, Sub Category, Repeat count, Remarks
$rows = @(
    @(42, "22-12-2025", "varun",  8289959250, "05-01-2026", "SHAIKRIZWAN",        "Loss", "PRODUCT",           "REQUIRED MODEL NOT AVAILABLE", "-", "NEED MINIMAL WORK"),
    @(43, "22-12-2025", "aparna", 8590441298, "23-12-2025", "MUHAMMED ROSHAN C V","Loss", "SIZE NOT SUITABLE", "SIZE TOO SMALL",               "-", "NEED BIG SIZE"),
    @(44, "22-12-2025", "ashwin", 9061864064, "01-02-2026", "MUHAMMED ROSHAN C V","Loss", "ENQUIRY",           "ENQUIRY WITHOUT BRIDE/FAMILY", "-", "WILL REVISIT"),
    @(45, "23-12-2025", "Akash",  9037331112, "03-01-2026", "MUHAMMED ROSHAN C V","Loss", "SIZE NOT SUITABLE", "SIZE TOO SMALL",               "-", "need bigger size IW"),
    @(46, "24-12-2025", "anuraj", 8289893408, "26-01-2026", "SHAIKRIZWAN",        "Loss", "ENQUIRY",           "ENQUIRY WITHOUT BRIDE/FAMILY", "-", "need to visit our vadakara store"),
    @(47, "25-12-2025", "arjun",  9544974956, "04-01-2026", "MUHAMMED ROSHAN C V","Loss", "ENQUIRY",           "ENQUIRY WITHOUT BRIDE/FAMILY", "-", "need to visit other store")
)

$startRow = 44
$lastExistingRow = $startRow - 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A (#) and D (Contact) are numeric - copy the number format from
    # the row above so the new cells share the existing style (s="65")
    # instead of minting a brand-new style entry.
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item($lastExistingRow, 1).NumberFormat

    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($lastExistingRow, 4).NumberFormat

    # Remaining columns are free-text. Some values look like ambiguous
    # dd-mm-yyyy dates (day <= 12), which Excel would otherwise auto-convert
    # to a date serial. Force text mode while assigning, then restore the
    # default "Normal" style so no stray number-format style lingers on the
    # cell (matching the plain, unstyled text cells used elsewhere).
    $textCols = @(2, 3, 5, 6, 7, 8, 9, 10, 11)
    foreach ($col in $textCols) {
        $ws.Cells.Item($r, $col).NumberFormat = "@"
        $ws.Cells.Item($r, $col).Value = $data[$col - 1]
        $ws.Cells.Item($r, $col).Style = "Normal"
    }
}
